$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'91.858.14"
$ws.Range("E2").Value = "'  -3.08%  "

# Row 3
$ws.Range("D3").Value = "'3.285.80"
$ws.Range("E3").Value = "'  -5.12%  "

# Row 4
$ws.Range("E4").Value = "'  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'227.21"
$ws.Range("E5").Value = "'  -5.31%  "

# Row 6
$ws.Range("D6").Value = "'608.02"
$ws.Range("E6").Value = "'  -5.70%  "

# Row 7
$ws.Range("D7").Value = "'1.35"
$ws.Range("E7").Value = "'  -7.83%  "

# Row 8
$ws.Range("D8").Value = "'0.377"
$ws.Range("E8").Value = "'  -6.88%  "

# Row 9
$ws.Range("E9").Value = "'  -0.04%  "

# Row 10
$ws.Range("D10").Value = "'0.927"
$ws.Range("E10").Value = "'  -8.37%  "

# Row 11
$ws.Range("D11").Value = "'3.283.81"
$ws.Range("E11").Value = "'  -5.12%  "

# Row 12
$ws.Range("D12").Value = "'41.37"
$ws.Range("E12").Value = "'  -1.40%  "

# Row 13
$ws.Range("D13").Value = "'0.191"
$ws.Range("E13").Value = "'  -4.01%  "

# Row 14
$ws.Range("D14").Value = "'5.90"
$ws.Range("E14").Value = "'  -4.00%  "

# Row 15
$ws.Range("D15").Value = "'91.676.31"
$ws.Range("E15").Value = "'  -3.16%  "

# Row 16
$ws.Range("D16").Value = "'3.896.17"
$ws.Range("E16").Value = "'  -5.25%  "

# Row 17
$ws.Range("D17").Value = "'0.0000239"
$ws.Range("E17").Value = "'  -6.65%  "

# Row 18
$ws.Range("D18").Value = "'7.97"
$ws.Range("E18").Value = "'  -6.66%  "

# Row 19
$ws.Range("D19").Value = "'3.283.13"
$ws.Range("E19").Value = "'  -5.40%  "

# Row 20
$ws.Range("D20").Value = "'16.97"
$ws.Range("E20").Value = "'  -5.29%  "

# Row 21
$ws.Range("D21").Value = "'10.62"
$ws.Range("E21").Value = "'  -6.92%  "

# Row 22
$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "'  +6.05%  "

# Row 23
$ws.Range("D23").Value = "'482.44"
$ws.Range("E23").Value = "'  -4.02%  "

# Row 24
$ws.Range("D24").Value = "'0.435"
$ws.Range("E24").Value = "'  -15.17%  "

# Row 25
$ws.Range("D25").Value = "'0.0000176"
$ws.Range("E25").Value = "'  -8.72%  "

# Row 26
$ws.Range("D26").Value = "'5.99"
$ws.Range("E26").Value = "'  -9.80%  "

# Row 27
$ws.Range("D27").Value = "'88.50"

# Row 28
$ws.Range("D28").Value = "'11.59"
$ws.Range("E28").Value = "'  -4.72%  "

# Row 29
$ws.Range("D29").Value = "'3.457.72"
$ws.Range("E29").Value = "'  -5.24%  "

# Row 30
$ws.Range("E30").Value = "'  +0.10%  "

# Row 31
$ws.Range("D31").Value = "'10.86"
$ws.Range("E31").Value = "'  -7.52%  "

# Row 32
$ws.Range("D32").Value = "'0.135"
$ws.Range("E32").Value = "'  -2.16%  "

# Row 33
$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.57"
$ws.Range("E33").Value = "'  -7.13%  "

# Row 34
$ws.Range("B34").Value = "'Binance-PegBSC-USD"
$ws.Range("C34").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.44%  "

# Row 35
$ws.Range("E35").Value = "'  -8.18%  "

# Row 36
$ws.Range("D36").Value = "'27.72"
$ws.Range("E36").Value = "'  -10.71%  "

# Row 37
$ws.Range("D37").Value = "'0.514"
$ws.Range("E37").Value = "'  -9.58%  "

# Row 38
$ws.Range("D38").Value = "'533.19"
$ws.Range("E38").Value = "'  +0.37%  "

# Row 39
$ws.Range("E39").Value = "'  -0.01%  "

# Row 40
$ws.Range("D40").Value = "'7.22"
$ws.Range("E40").Value = "'  -7.07%  "

# Row 41
$ws.Range("E41").Value = "'  -3.57%  "

# Row 42
$ws.Range("D42").Value = "'1.34"
$ws.Range("E42").Value = "'  -7.40%  "

# Row 43
$ws.Range("D43").Value = "'0.846"
$ws.Range("E43").Value = "'  -8.77%  "

# Row 44
$ws.Range("E44").Value = "'  -1.38%  "

# Row 45
$ws.Range("B45").Value = "'MantraDAO"
$ws.Range("C45").Value = "'https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.56"
$ws.Range("E45").Value = "'  +1.39%  "

# Row 46
$ws.Range("B46").Value = "'ImmutableX"
$ws.Range("C46").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.65"
$ws.Range("E46").Value = "'  -3.13%  "

# Row 47
$ws.Range("E47").Value = "'  -4.08%  "

# Row 48
$ws.Range("D48").Value = "'5.24"
$ws.Range("E48").Value = "'  -8.42%  "

# Row 49
$ws.Range("B49").Value = "'Stacks"
$ws.Range("C49").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.05"
$ws.Range("E49").Value = "'  -4.69%  "

# Row 50
$ws.Range("B50").Value = "'OKB"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'51.25"
$ws.Range("E50").Value = "'  -3.87%  "

# Row 51
$ws.Range("D51").Value = "'7.80"
$ws.Range("E51").Value = "'  -3.60%  "

